$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Початок 43 поля" note to "Початок 42 поля" for the first four
# burial records (column E, rows 2-5).
$ws.Range("E2").Value = "Початок 42 поля"
$ws.Range("E3").Value = "Початок 42 поля"
$ws.Range("E4").Value = "Початок 42 поля"
$ws.Range("E5").Value = "Початок 42 поля"

# Row 24 no longer needs its taller custom height - let Excel re-fit it to
# the default row height.
$ws.Rows.Item(24).AutoFit()

# Reflect the scrolled/selected view state left behind by the edit.
$ws.Range("A57").Select()
$excel.ActiveWindow.ScrollRow = 41
